# Generate Report for Handback
#
# Two e2e test-fixture files were renamed/regenerated:
#   a70c4802-e0a0-44ae-9b0f-23258a963346.md -> 4c4da03c-d97d-4d02-ac62-9daaa187695b.md
#   d5427d45-234e-468a-975d-98b7678503d4.md -> ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md
# and the handback pipeline re-ran, producing new xliff names / timestamps.
# This updates the three report sheets (Overview, zh-cn, de-de) accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("B2").Value = "e2e\4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("G2").Value = "2016-08-30 09:32:11"

$ws.Range("A3").Value = "ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("B3").Value = "e2e\ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("G3").Value = "2016-08-30 09:32:11"

# ---------------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("G2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-30 09:31:58"
$ws.Range("I2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("J2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-30 09:34:23"

$ws.Range("A3").Value = "ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("G3").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-30 09:31:58"
$ws.Range("I3").Value = "ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("J3").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-30 09:34:23"

# ---------------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("G2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.de-de.xlf"
$ws.Range("H2").Value = "2016-08-30 09:32:11"
$ws.Range("I2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$ws.Range("J2").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.de-de.xlf"
$ws.Range("K2").Value = "2016-08-30 09:34:31"

$ws.Range("A3").Value = "ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("G3").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.de-de.xlf"
$ws.Range("H3").Value = "2016-08-30 09:32:11"
$ws.Range("I3").Value = "ffff95f47282-4fd1-4cec-966d-73434d08c4d3.md"
$ws.Range("J3").Value = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.de-de.xlf"
$ws.Range("K3").Value = "2016-08-30 09:34:31"

Write-Output "handback-status report updated"
